$wb = $excel.ActiveWorkbook

# --- Duplicate "Level2_Econ" placing the copy right after it -------------
# This keeps the original sheet (sheetId=2) in position 1 and gives the new
# copy a brand-new sheetId, matching the workbook.xml sheet list in the diff:
#   Level2_EconTest (sheetId=2, was "Level2_Econ")
#   Level2_Econ     (sheetId=7, new copy)
#   Level1, Level2_Technical, Level2_Infrastructure, Level2_serviceStandard
$wsEconOriginal = $wb.Worksheets.Item("Level2_Econ")
$wsEconOriginal.Copy($null, $wsEconOriginal) | Out-Null

$wsTest = $wb.Worksheets.Item("Level2_Econ")
$wsTest.Name = "Level2_EconTest"

$wsEcon = $wb.Worksheets.Item("Level2_Econ (2)")
$wsEcon.Name = "Level2_Econ"

# --- Add two "empty" placeholder criteria columns on the test sheet ------
# (AHP needs >= 3 criteria; the original sheet only had 1 pairwise-comparison
# column, so two "empty" columns with neutral value 0 are appended.)
$wsTest.Range("B1").Value = "empty"
$wsTest.Range("C1").Value = "empty"
$wsTest.Range("B2").Value = 0
$wsTest.Range("C2").Value = 0

# --- Fix up selections on a couple of other sheets ------------------------
$wsTech = $wb.Worksheets.Item("Level2_Technical")
$wsTech.Range("A2").Select() | Out-Null

$wsTest.Range("A1:A2").Select() | Out-Null

# Make "Level2_Econ" (the fresh copy) the active / selected tab on save,
# matching the unchanged activeTab index in bookViews.
$wsEcon.Activate()
